$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "To" time for the 4th data row (row 5) from 21:00 to 23:00.
$ws.Range("C5").Value = 0.958333333333333

# Move the selection/active cell to C6 (matches author's final cursor position).
$ws.Range("C6").Select()
